$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Investor name) - rows 5-8
$ws.Range("A5").Value = "Investor 3"
$ws.Range("A6").Value = "Investor 4"
$ws.Range("A7").Value = "Investor 5"
$ws.Range("A8").Value = "Investor 6"

# Column B (First Name) - rows 5-8
$ws.Range("B5").Value = "Emp3"
$ws.Range("B6").Value = "Emp4"
$ws.Range("B7").Value = "Emp5"
$ws.Range("B8").Value = "Emp6"

# Column C (Last Name) - rows 5-8
$ws.Range("C5").Value = "L3"
$ws.Range("C6").Value = "L4"
$ws.Range("C7").Value = "L5"
$ws.Range("C8").Value = "L6"

# Column D (Email, hyperlinked) - rows 5-8
$ws.Range("D5").Value = "emp1@investor3.com"
$ws.Range("D6").Value = "emp1@investor4.com"
$ws.Range("D7").Value = "emp1@investor5.com"
$ws.Range("D8").Value = "emp1@investor6.com"

# Remaining columns - row 5
$ws.Range("E5").Value = "Yes"
$ws.Range("G5").Value = "IN(91)"
$ws.Range("H5").Value = 1111111111
$ws.Range("I5").Value = "Yes"
$ws.Range("J5").Value = "Yes"
$ws.Range("K5").Value = "No"

# Remaining columns - row 6
$ws.Range("E6").Value = "Yes"
$ws.Range("G6").Value = "IN(91)"
$ws.Range("H6").Value = 2222222222
$ws.Range("I6").Value = "Yes"
$ws.Range("J6").Value = "Yes"
$ws.Range("K6").Value = "No"

# Remaining columns - row 7
$ws.Range("E7").Value = "Yes"
$ws.Range("G7").Value = "IN(91)"
$ws.Range("H7").Value = 3333333333
$ws.Range("I7").Value = "Yes"
$ws.Range("J7").Value = "Yes"
$ws.Range("K7").Value = "No"

# Remaining columns - row 8 (partial row - no H, I, J, K)
$ws.Range("E8").Value = "Yes"
$ws.Range("G8").Value = "IN(91)"

# Add the mailto hyperlinks (creates relationship + hyperlink element)
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:emp1@investor3.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:emp1@investor4.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:emp1@investor5.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:emp1@investor6.com") | Out-Null

# Restore the Hyperlink cell style (Hyperlinks.Add applies its own style variant)
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"
$ws.Range("D7").Style = "Hyperlink"
$ws.Range("D8").Style = "Hyperlink"

# Update selection to match final cursor position
$ws.Range("K8").Select() | Out-Null

Write-Output "done"
